$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-20 00:50:20"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-20 00:50:15"
$wsZhCn.Range("K4").Value = "2016-08-20 00:50:34"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-20 00:50:20"
$wsDeDe.Range("K4").Value = "2016-08-20 00:50:48"
